# Add team record (Wins / Losses / Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled the same as the rest of row 1 (bold, bordered, centered).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill every data row (2 through 50) with the team's win/loss/tie record.
$ws.Range("AD2:AD50").Value = 76
$ws.Range("AE2:AE50").Value = 86
$ws.Range("AF2:AF50").Value = 0
